$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 3 and row 4 for columns D, J, K, L, M, P
# (Fecha, Volumen, Precio mínimo, Precio máximo, Precio promedio ponderado, Precio $/Kg)

$ws.Range("D3").Value = 44827
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 31000
$ws.Range("M3").Value = 30500
$ws.Range("P3").Value = 1220

$ws.Range("D4").Value = 44414
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 31000
$ws.Range("L4").Value = 32000
$ws.Range("M4").Value = 31500
$ws.Range("P4").Value = 1260
